$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DatosCP")

# --- Fix border style on C4/C5 (reuse the xf already used by A2/A3: font underline + thin box border) ---
$ws.Range("A2").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- New rows 6-9: propagate the row-5 formatting first ---
$ws.Range("A5:D5").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)
$ws.Range("A5:D5").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)
$ws.Range("A5:D5").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)
$ws.Range("A5:D5").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Fill in the test case data (order matches the authoring sequence) ---
$ws.Range("A6").Value = "CP005_checkThatProductAppearsonWantedList"
$ws.Range("A7").Value = "CP006_searchProductNotAddedToWishList"
$ws.Range("A5").Value = "CP004_addProductToWishlist"
$ws.Range("D8").Value = "Product successfully added to your shopping cart"
$ws.Range("A8").Value = "CP007_agregarProductoSinLogin"
$ws.Range("D9").Value = "Your shopping cart is empty."
$ws.Range("A9").Value = "CP008_emptyCartMessage"

# --- Login credentials reused on every new row ---
$ws.Range("B6").Value = "mailtestautomation001@gmail.com"
$ws.Range("C6").Value = "admin123"
$ws.Range("B7").Value = "mailtestautomation001@gmail.com"
$ws.Range("C7").Value = "admin123"
$ws.Range("B8").Value = "mailtestautomation001@gmail.com"
$ws.Range("C8").Value = "admin123"
$ws.Range("B9").Value = "mailtestautomation001@gmail.com"
$ws.Range("C9").Value = "admin123"

# --- Add hyperlinks on column B for the new rows (same mailto target as rows 4-5) ---
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:mailtestautomation001@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:mailtestautomation001@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:mailtestautomation001@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:mailtestautomation001@gmail.com")

# Hyperlinks.Add() stamps its own style - reapply the clean row-5 style afterwards
$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Final selection ---
$ws.Range("C9").Select()
